# Applies the "Updated symbol list" crypto-price refresh described by the
# commit diff: cells D2:E51 (Price / Volume(1h)) get refreshed numeric
# strings for the rows whose coin still trades (rows with "--" placeholders
# and a few untouched rows are intentionally left alone).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    # Force text interpretation so values like "324.72" / "-1.78%" are
    # stored as literal strings (matching the sheet's existing inline-string
    # convention) instead of being auto-coerced into numbers/percentages.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    # Reset to the default "Normal" style afterwards so we don't leave a
    # stray number-format style attached to the cell (the source cells carry
    # no explicit style).
    $cell.Style = "Normal"
}

Set-TextCell 2 4 "324.72"
Set-TextCell 2 5 "-1.78%"
Set-TextCell 3 4 "39.36"
Set-TextCell 3 5 "-2.74%"
Set-TextCell 4 4 "5.699"
Set-TextCell 4 5 "7.71%"
Set-TextCell 5 4 "0.08006"
Set-TextCell 5 5 "-1.17%"
Set-TextCell 6 4 "2.000"
Set-TextCell 6 5 "4.22%"
Set-TextCell 7 4 "4.495"
Set-TextCell 7 5 "-0.41%"
Set-TextCell 8 4 "8.594"
Set-TextCell 8 5 "-0.22%"
Set-TextCell 9 4 "2.974"
Set-TextCell 9 5 "-0.12%"
Set-TextCell 10 4 "0.9258"
Set-TextCell 10 5 "-1.09%"
Set-TextCell 11 4 "0.1237"
Set-TextCell 12 4 "0.1976"
Set-TextCell 12 5 "-0.24%"
Set-TextCell 13 4 "8.718"
Set-TextCell 13 5 "21.51%"
Set-TextCell 14 4 "0.09161"
Set-TextCell 14 5 "-1.09%"
Set-TextCell 15 5 "1.15%"
Set-TextCell 16 4 "0.1048"
Set-TextCell 16 5 "9.53%"
Set-TextCell 17 4 "0.001294"
Set-TextCell 17 5 "-2.75%"
Set-TextCell 18 4 "0.006135"
Set-TextCell 18 5 "-1.20%"
Set-TextCell 19 5 "-0.78%"
Set-TextCell 21 4 "0.1371"
Set-TextCell 21 5 "3.72%"
Set-TextCell 22 5 "-5.76%"
Set-TextCell 23 4 "0.04407"
Set-TextCell 23 5 "-0.21%"
Set-TextCell 24 4 "0.001260"
Set-TextCell 24 5 "3.38%"
Set-TextCell 25 4 "0.004630"
Set-TextCell 25 5 "7.91%"
Set-TextCell 26 4 "0.0001150"
Set-TextCell 26 5 "-3.38%"
Set-TextCell 39 4 "0.02490"
Set-TextCell 39 5 "-0.48%"
Set-TextCell 40 4 "0.05335"
Set-TextCell 40 5 "2.77%"
Set-TextCell 41 4 "0.007462"
Set-TextCell 41 5 "-2.95%"
Set-TextCell 42 4 "0.009621"
Set-TextCell 42 5 "3.64%"
Set-TextCell 43 5 "-1.40%"
Set-TextCell 44 4 "0.002116"
Set-TextCell 44 5 "-2.52%"
Set-TextCell 45 4 "0.01033"
Set-TextCell 45 5 "-0.43%"
Set-TextCell 46 4 "0.00006716"
Set-TextCell 46 5 "1.10%"
Set-TextCell 47 4 "0.00000000750"
Set-TextCell 47 5 "-0.02%"
Set-TextCell 48 4 "0.002970"
Set-TextCell 48 5 "-11.20%"
Set-TextCell 49 4 "0.002291"
Set-TextCell 49 5 "-4.56%"
Set-TextCell 50 4 "0.00002100"
Set-TextCell 50 5 "-0.02%"
Set-TextCell 51 4 "0.0002000"
Set-TextCell 51 5 "-0.02%"
